# Update Fgf15-Fgfr1 NATMI sheet with recomputed TPM-based values.
# Sending cluster column (A) and the dependent ligand/receptor/edge
# statistic columns (E:T) change for every data row; columns B, C, D,
# K and L (ligand/receptor symbol, target cluster, receptor cell counts)
# are unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05008433333333333
$ws.Range("H2").Value = 0.150253
$ws.Range("I2").Value = 0.3054767171413236
$ws.Range("J2").Value = 0.3054767171413236
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.48767733333333
$ws.Range("N2").Value = 31.463032
$ws.Range("O2").Value = 0.1222087640673552
$ws.Range("P2").Value = 0.1222087640673552
$ws.Range("Q2").Value = 0.5252683274551111
$ws.Range("R2").Value = 4.727414947095999
$ws.Range("S2").Value = 0.03733193205319423
$ws.Range("T2").Value = 0.03733193205319423

$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05008433333333333
$ws.Range("H3").Value = 0.150253
$ws.Range("I3").Value = 0.3054767171413236
$ws.Range("J3").Value = 0.3054767171413236
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 62.99699166666667
$ws.Range("N3").Value = 188.990975
$ws.Range("O3").Value = 0.7340790765058636
$ws.Range("P3").Value = 0.7340790765058635
$ws.Range("Q3").Value = 3.155162329630556
$ws.Range("R3").Value = 28.396460966675
$ws.Range("S3").Value = 0.2242440664131458
$ws.Range("T3").Value = 0.2242440664131457

$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05008433333333333
$ws.Range("H4").Value = 0.150253
$ws.Range("I4").Value = 0.3054767171413236
$ws.Range("J4").Value = 0.3054767171413236
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3322793333333333
$ws.Range("N4").Value = 0.996838
$ws.Range("O4").Value = 0.003871919907635547
$ws.Range("P4").Value = 0.003871919907635547
$ws.Range("Q4").Value = 0.01664198889044444
$ws.Range("R4").Value = 0.149777900014
$ws.Range("S4").Value = 0.001182781382418644
$ws.Range("T4").Value = 0.001182781382418644

$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05008433333333333
$ws.Range("H5").Value = 0.150253
$ws.Range("I5").Value = 0.3054767171413236
$ws.Range("J5").Value = 0.3054767171413236
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.25458433333333
$ws.Range("N5").Value = 30.763753
$ws.Range("O5").Value = 0.1194926233493133
$ws.Range("P5").Value = 0.1194926233493133
$ws.Range("Q5").Value = 0.5135940199454445
$ws.Range("R5").Value = 4.622346179509
$ws.Range("S5").Value = 0.03650221430335291
$ws.Range("T5").Value = 0.03650221430335291

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05008433333333333
$ws.Range("H6").Value = 0.150253
$ws.Range("I6").Value = 0.3054767171413236
$ws.Range("J6").Value = 0.3054767171413236
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7572163333333334
$ws.Range("N6").Value = 2.271649
$ws.Range("O6").Value = 0.008823543029319092
$ws.Range("P6").Value = 0.00882354302931909
$ws.Range("Q6").Value = 0.03792467524411111
$ws.Range("R6").Value = 0.341322077197
$ws.Range("S6").Value = 0.002695386958151606
$ws.Range("T6").Value = 0.002695386958151606

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05008433333333333
$ws.Range("H7").Value = 0.150253
$ws.Range("I7").Value = 0.3054767171413236
$ws.Range("J7").Value = 0.3054767171413236
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9889696666666666
$ws.Range("N7").Value = 2.966909
$ws.Range("O7").Value = 0.01152407314051338
$ws.Range("P7").Value = 0.01152407314051338
$ws.Range("Q7").Value = 0.04953188644188889
$ws.Range("R7").Value = 0.445786977977
$ws.Range("S7").Value = 0.00352033603106053
$ws.Range("T7").Value = 0.00352033603106053

$ws.Range("A8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1138703333333333
$ws.Range("H8").Value = 0.341611
$ws.Range("I8").Value = 0.6945232828586764
$ws.Range("J8").Value = 0.6945232828586764
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.48767733333333
$ws.Range("N8").Value = 31.463032
$ws.Range("O8").Value = 0.1222087640673552
$ws.Range("P8").Value = 0.1222087640673552
$ws.Range("Q8").Value = 1.194235313839111
$ws.Range("R8").Value = 10.748117824552
$ws.Range("S8").Value = 0.084876832014161
$ws.Range("T8").Value = 0.084876832014161

$ws.Range("A9").Value = "MuSCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1138703333333333
$ws.Range("H9").Value = 0.341611
$ws.Range("I9").Value = 0.6945232828586764
$ws.Range("J9").Value = 0.6945232828586764
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 62.99699166666667
$ws.Range("N9").Value = 188.990975
$ws.Range("O9").Value = 0.7340790765058636
$ws.Range("P9").Value = 0.7340790765058635
$ws.Range("Q9").Value = 7.173488440080556
$ws.Range("R9").Value = 64.56139596072499
$ws.Range("S9").Value = 0.5098350100927178
$ws.Range("T9").Value = 0.5098350100927177

$ws.Range("A10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1138703333333333
$ws.Range("H10").Value = 0.341611
$ws.Range("I10").Value = 0.6945232828586764
$ws.Range("J10").Value = 0.6945232828586764
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3322793333333333
$ws.Range("N10").Value = 0.996838
$ws.Range("O10").Value = 0.003871919907635547
$ws.Range("P10").Value = 0.003871919907635547
$ws.Range("Q10").Value = 0.03783675844644444
$ws.Range("R10").Value = 0.340530826018
$ws.Range("S10").Value = 0.002689138525216903
$ws.Range("T10").Value = 0.002689138525216903

$ws.Range("A11").Value = "MuSCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1138703333333333
$ws.Range("H11").Value = 0.341611
$ws.Range("I11").Value = 0.6945232828586764
$ws.Range("J11").Value = 0.6945232828586764
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.25458433333333
$ws.Range("N11").Value = 30.763753
$ws.Range("O11").Value = 0.1194926233493133
$ws.Range("P11").Value = 0.1194926233493133
$ws.Range("Q11").Value = 1.167692936231445
$ws.Range("R11").Value = 10.509236426083
$ws.Range("S11").Value = 0.08299040904596042
$ws.Range("T11").Value = 0.08299040904596042

$ws.Range("A12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1138703333333333
$ws.Range("H12").Value = 0.341611
$ws.Range("I12").Value = 0.6945232828586764
$ws.Range("J12").Value = 0.6945232828586764
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.7572163333333334
$ws.Range("N12").Value = 2.271649
$ws.Range("O12").Value = 0.008823543029319092
$ws.Range("P12").Value = 0.00882354302931909
$ws.Range("Q12").Value = 0.08622447628211112
$ws.Range("R12").Value = 0.776020286539
$ws.Range("S12").Value = 0.006128156071167486
$ws.Range("T12").Value = 0.006128156071167486

$ws.Range("A13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1138703333333333
$ws.Range("H13").Value = 0.341611
$ws.Range("I13").Value = 0.6945232828586764
$ws.Range("J13").Value = 0.6945232828586764
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.9889696666666666
$ws.Range("N13").Value = 2.966909
$ws.Range("O13").Value = 0.01152407314051338
$ws.Range("P13").Value = 0.01152407314051338
$ws.Range("Q13").Value = 0.1126143055998889
$ws.Range("R13").Value = 1.013528750399
$ws.Range("S13").Value = 0.008003737109452848
$ws.Range("T13").Value = 0.008003737109452846

